# Weekly update: insert a new daily price record as a new row 366 in the
# "Agrícola del Norte S.A. de Arica - Plátano" price log, pushing the
# existing rows 366-409 down to 367-410 (last row's data lands in the new
# row 410). This mirrors Excel's "Insert Row" (shift cells down) behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 366, shifting rows 366:409 down to 367:410.
$ws.Rows(366).Insert()

# Populate the newly inserted row 366 with the new record.
$ws.Range("A366").Value = 1
$ws.Range("B366").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C366").Value = "Arica y Parinacota"
$ws.Range("D366").Value = 45127
$ws.Range("E366").Value = 15
$ws.Range("F366").Value = "Fruta"
$ws.Range("G366").Value = 100108
$ws.Range("H366").Value = "Tropicales y subtropicales"
$ws.Range("I366").Value = 100108006
$ws.Range("J366").Value = "Plátano"
$ws.Range("K366").Value = "Sin especificar"
$ws.Range("L366").Value = "Pintón"
$ws.Range("M366").Value = 250
$ws.Range("N366").Value = 19000
$ws.Range("O366").Value = 20000
$ws.Range("P366").Value = 19400
$ws.Range("Q366").Value = "`$/caja 20 kilos"
$ws.Range("R366").Value = "Ecuador"
$ws.Range("S366").Value = 970
$ws.Range("T366").Value = 20
